# TC_51822: update the Alarm Load test value for the Generic Sounder row
# and leave the sheet's active selection on the next cell to verify (E10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Devices")

# Alarm Load for the "Generic Sounder" device (row 11) changes from 0 to 25
$ws.Range("F11").Value = 25

# Move/leave the active selection at E10, matching the saved view state
$ws.Activate()
$ws.Range("E10").Select()
